$d = $word.ActiveDocument

# --- First paragraph: update the ID placeholder text and paragraph formatting ---
$p = $d.Paragraphs(1)

# Replace the paragraph's text (both runs) with the new single-run text,
# excluding the trailing paragraph mark. This collapses the two runs
# (the ID text run + the trailing " " run) down to a single run, dropping
# the now-unneeded trailing space.
$r = $p.Range
$null = $r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_SUBPART_5342_71__ID**"

# Add a paragraph border (top/left/bottom/right) with 5pt spacing, and
# change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$pf = $p.Format
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25
